$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStrings = @(
    "SCRIPT/P02P01A/um1205.ssb",
    "SCRIPT/P02P01A/um1301.ssb",
    "SCRIPT/P02P01A/um1304.ssb",
    "SCRIPT/P02P01A/um1307.ssb",
    "SCRIPT/P02P01A/um1402.ssb",
    "SCRIPT/P02P01A/um1405.ssb",
    "SCRIPT/P02P01A/um1502.ssb",
    "SCRIPT/P02P01A/um1602.ssb",
    "SCRIPT/P02P01A/um1605.ssb"
)

$templateRow = $ws.Rows.Item(7)
$templateHeight = $templateRow.RowHeight

$row = 8
foreach ($s in $newStrings) {
    $templateRow.Copy()
    $ws.Rows.Item($row).PasteSpecial()
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $s
    $ws.Rows.Item($row).RowHeight = $templateHeight
    $row = $row + 1
}

$ws.Range("C16").Select()
